$d = $word.ActiveDocument

# Locate the end of the last bibliography paragraph
# ("...o século 21. Editora Campus. 2010.").
$endOfBiblio = $d.Content
$endOfBiblio.Find.Execute("o século 21. Editora Campus. 2010.", $true, $false, $false, `
                           $false, $false, $true, 1, $false, "", 0) | Out-Null
$endOfBiblio.Expand(4) | Out-Null

# Locate the paragraph with the site footer/copyright text that
# (together with the two blank/page-break paragraphs right after the
# bibliography) is being removed.
$footer = $d.Content
$footer.Find.Execute("© 2020", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$footer.Expand(4) | Out-Null

# Remove everything from right after the bibliography paragraph through
# the end (incl. paragraph mark) of the footer paragraph. This deletes
# the blank paragraph, the blank page-break paragraph, and the footer
# paragraph, leaving the paragraphs that follow untouched.
$d.Range($endOfBiblio.End, $footer.End).Delete() | Out-Null
